$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.2301287744451176
$ws.Range("E3").Value = 0.2531497705465164
$ws.Range("F3").Value = 0.3792262962517786
$ws.Range("D4").Value = 0.622085491896273
$ws.Range("E4").Value = 0.715803612783696
$ws.Range("F4").Value = 1.410839771178848
$ws.Range("D5").Value = 0.231404064342024
$ws.Range("E5").Value = 0.3186326968857948
$ws.Range("F5").Value = 0.4753388806427217
$ws.Range("D7").Value = 0.1472510254554347
$ws.Range("E7").Value = 0.2140058763679858
$ws.Range("F7").Value = 0.2750689986460114
$ws.Range("D8").Value = 0.5256880469229492
$ws.Range("E8").Value = 1.098314156725016
$ws.Range("F8").Value = 1.311768761073552
$ws.Range("D9").Value = 0.340718843625071
$ws.Range("E9").Value = 0.5875718223060855
$ws.Range("F9").Value = 0.941394689181267
$ws.Range("D10").Value = 0.1619387528663363
$ws.Range("E10").Value = 0.1834216150915343
$ws.Range("F10").Value = 0.2842812774501192
$ws.Range("D12").Value = 1.063454797558633
$ws.Range("E12").Value = 1.063454797558633
$ws.Range("F12").Value = 1.063454797558633
$ws.Range("D13").Value = 0.08797278745859347
$ws.Range("E13").Value = 0.1388735943783553
$ws.Range("F13").Value = 0.1709720771502557
$ws.Range("D16").Value = 0.7984837159027873
$ws.Range("E16").Value = 0.8966322525315829
$ws.Range("F16").Value = 1.013612461361046
$ws.Range("D17").Value = 0.2591769352772525
$ws.Range("E17").Value = 0.286278606783992
$ws.Range("F17").Value = 0.3779984695809697
$ws.Range("D19").Value = 0.2490244497601257
$ws.Range("E19").Value = 0.2868897572640994
$ws.Range("F19").Value = 0.3931129945572671
$ws.Range("D20").Value = 0.594910893707619
$ws.Range("E20").Value = 0.749578755244962
$ws.Range("F20").Value = 1.109549047795467
$ws.Range("D24").Value = 0.7653095040833233
$ws.Range("E24").Value = 0.9643827264910279
$ws.Range("F24").Value = 1.005395117494798
$ws.Range("D25").Value = 0.04385119017408517
$ws.Range("E25").Value = 0.06274509803921569
$ws.Range("F25").Value = 0.1058665827941529
$ws.Range("D27").Value = 0.6271560677089354
$ws.Range("E27").Value = 0.8731027926322045
$ws.Range("F27").Value = 2.049484120224041
$ws.Range("D28").Value = 1.292880440132429
$ws.Range("E28").Value = 1.765034646914061
$ws.Range("F28").Value = 2.250614141699732
$ws.Range("D29").Value = 0.224719757213062
$ws.Range("E29").Value = 0.335096784774885
$ws.Range("F29").Value = 0.5730657367164054
$ws.Range("D30").Value = 0.4621151927714648
$ws.Range("E30").Value = 0.5435085557188226
$ws.Range("F30").Value = 0.8551305927776516
$ws.Range("D31").Value = 0.1077957545940685
$ws.Range("E31").Value = 0.1641717369349599
$ws.Range("F31").Value = 0.200518436383942
$ws.Range("D32").Value = 0.7243747686308436
$ws.Range("E32").Value = 1.044690644282478
$ws.Range("F32").Value = 1.2761306200802
$ws.Range("D34").Value = 0.1521371538257352
$ws.Range("E34").Value = 0.2470020714508669
$ws.Range("F34").Value = 0.3966944982048929
$ws.Range("D35").Value = 7.866091611315241
$ws.Range("E35").Value = 12.56559700431341
$ws.Range("F35").Value = 18.37973599791526
$ws.Range("D36").Value = 1.093835496047211
$ws.Range("E36").Value = 1.688259240258801
$ws.Range("F36").Value = 2.072730515421942
$ws.Range("D37").Value = 0.1341224666758354
$ws.Range("E37").Value = 0.2126072592625441
$ws.Range("F37").Value = 0.4800733921942756
$ws.Range("D40").Value = 1.253931337924155
$ws.Range("E40").Value = 1.613871382538543
$ws.Range("F40").Value = 2.159850156652101
$ws.Range("D41").Value = 0.4055990638889103
$ws.Range("E41").Value = 0.4183592204963659
$ws.Range("F41").Value = 0.5239198976767825
$ws.Range("D43").Value = 0.310672514619883
$ws.Range("E43").Value = 0.4227411027106767
$ws.Range("F43").Value = 0.4528457268346245
$ws.Range("D44").Value = 0.5067555239035585
$ws.Range("E44").Value = 0.915501537675697
$ws.Range("F44").Value = 0.9923590603841175
$ws.Range("D45").Value = 0.4935939800344355
$ws.Range("E45").Value = 0.7104285125446085
$ws.Range("F45").Value = 1.061763089777439
$ws.Range("D46").Value = 0.2230930592240576
$ws.Range("E46").Value = 0.2878711045546745
$ws.Range("F46").Value = 0.5535336562265663
$ws.Range("D48").Value = 1.228065015479876
$ws.Range("E48").Value = 1.869227666128795
$ws.Range("F48").Value = 2.798037361714649
$ws.Range("D49").Value = 0.1311080233667974
$ws.Range("E49").Value = 0.1501780155422384
$ws.Range("F49").Value = 0.2143819258628985
